# Remove the "Ans:" code sample block (stack-validation script) that
# followed the "[1+2(ab] - Invalid Expression..." example paragraph,
# while keeping that paragraph and the trailing blank paragraph intact.

$d = $word.ActiveDocument

$startText = "Ans:"
$endText = "  print(""invalid"")"

$startPara = $null
$endPara = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq $startText -and $startPara -eq $null) {
        $startPara = $p
    }
    if ($t -eq $endText) {
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rng.Delete()
}
